# Applies the numeric data refresh captured by the commit diff across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR price-tracking sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1213.6666
$ws.Range("I98").Value = 594.8125
$ws.Range("J98").Value = 6164.5
$ws.Range("K98").Value = 594.8125
$ws.Range("L98").Value = 6164.5
$ws.Range("M98").Value = 903.1875
$ws.Range("N98").Value = -9160.5

$ws.Range("H122").Value = 1213.6666
$ws.Range("I122").Value = 594.8125
$ws.Range("J122").Value = 6164.5
$ws.Range("K122").Value = 1784.4375
$ws.Range("L122").Value = 18493.5
$ws.Range("M122").Value = 665.5625
$ws.Range("N122").Value = -23393.5

$ws.Range("H138").Value = 4606.8
$ws.Range("I138").Value = 1979.3
$ws.Range("J138").Value = 5657.8
$ws.Range("K138").Value = 5937.9
$ws.Range("L138").Value = 16973.4
$ws.Range("M138").Value = -797.8999999999996
$ws.Range("N138").Value = -27253.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7594.3237
$ws.Range("I32").Value = 2200.3462
$ws.Range("J32").Value = 25124.75
$ws.Range("K32").Value = 2200.3462
$ws.Range("L32").Value = 25124.75
$ws.Range("M32").Value = -1913.3462
$ws.Range("N32").Value = -25698.75

$ws.Range("H61").Value = 4320.1025
$ws.Range("I61").Value = 2916.9092
$ws.Range("J61").Value = 12037.667
$ws.Range("K61").Value = 2916.9092
$ws.Range("L61").Value = 12037.667
$ws.Range("M61").Value = -2704.9092
$ws.Range("N61").Value = -12461.667

$ws.Range("H74").Value = 2352.818
$ws.Range("I74").Value = 2352.818
$ws.Range("K74").Value = 2352.818
$ws.Range("M74").Value = -1478.818

$ws.Range("H77").Value = 2352.818
$ws.Range("I77").Value = 2352.818
$ws.Range("K77").Value = 11764.09
$ws.Range("M77").Value = -7396.09

$ws.Range("H122").Value = 5906.476
$ws.Range("I122").Value = 4930.857
$ws.Range("J122").Value = 7857.7144
$ws.Range("K122").Value = 14792.571
$ws.Range("L122").Value = 23573.1432
$ws.Range("M122").Value = -12342.571
$ws.Range("N122").Value = -28473.1432

$ws.Range("H132").Value = 4048.8
$ws.Range("I132").Value = 3061
$ws.Range("K132").Value = 9183
$ws.Range("M132").Value = -6653

$ws.Range("H136").Value = 4320.1025
$ws.Range("I136").Value = 2916.9092
$ws.Range("J136").Value = 12037.667
$ws.Range("K136").Value = 8750.7276
$ws.Range("L136").Value = 36113.001
$ws.Range("M136").Value = -6200.7276
$ws.Range("N136").Value = -41213.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3479.8
$ws.Range("I20").Value = 3399.8572
$ws.Range("K20").Value = 3399.8572
$ws.Range("M20").Value = -3152.8572

$ws.Range("H52").Value = 45000
$ws.Range("J52").Value = 45000
$ws.Range("L52").Value = 45000
$ws.Range("N52").Value = -45526

$ws.Range("H88").Value = 25000
$ws.Range("J88").Value = 25000
$ws.Range("L88").Value = 25000
$ws.Range("N88").Value = -25812

$ws.Range("H91").Value = 25000
$ws.Range("J91").Value = 25000
$ws.Range("L91").Value = 25000
$ws.Range("N91").Value = -27808

$ws.Range("H107").Value = 1530
$ws.Range("I107").Value = 1390.1
$ws.Range("K107").Value = 1390.1
$ws.Range("M107").Value = 529.9000000000001

$ws.Range("H112").Value = 49900
$ws.Range("J112").Value = 49900
$ws.Range("L112").Value = 49900
$ws.Range("N112").Value = -52854

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = $null

$ws.Range("H121").Value = 45000
$ws.Range("J121").Value = 45000
$ws.Range("L121").Value = 45000
$ws.Range("N121").Value = -48494

$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").Value = $null

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").Value = $null

$ws.Range("H134").Value = 5262.6665
$ws.Range("I134").Value = 5262.6665
$ws.Range("K134").Value = 15787.9995
$ws.Range("M134").Value = -13252.9995

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1568.44
$ws.Range("I58").Value = 1309.3684
$ws.Range("K58").Value = 1309.3684
$ws.Range("M58").Value = -1106.3684

$ws.Range("H122").Value = 102735.53
$ws.Range("J122").Value = 2267.3076
$ws.Range("L122").Value = 6801.9228
$ws.Range("N122").Value = -11701.9228

$ws.Range("H131").Value = 44540.363
$ws.Range("J131").Value = 48883.445
$ws.Range("L131").Value = 48883.445
$ws.Range("N131").Value = -58963.445

$ws.Range("H136").Value = 1568.44
$ws.Range("I136").Value = 1309.3684
$ws.Range("K136").Value = 3928.1052
$ws.Range("M136").Value = -1378.1052

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2857.4285
$ws.Range("I70").Value = 3001.4
$ws.Range("J70").Value = 2497.5
$ws.Range("K70").Value = 9004.200000000001
$ws.Range("L70").Value = 7492.5
$ws.Range("M70").Value = -8689.200000000001
$ws.Range("N70").Value = -8122.5

$ws.Range("H73").Value = 2857.4285
$ws.Range("I73").Value = 3001.4
$ws.Range("J73").Value = 2497.5
$ws.Range("K73").Value = 9004.200000000001
$ws.Range("L73").Value = 7492.5
$ws.Range("M73").Value = -7912.200000000001
$ws.Range("N73").Value = -9676.5

$ws.Range("H86").Value = 708.1818
$ws.Range("I86").Value = 656.8
$ws.Range("J86").Value = 751
$ws.Range("K86").Value = 1970.4
$ws.Range("L86").Value = 2253
$ws.Range("M86").Value = -784.3999999999999
$ws.Range("N86").Value = -4625

$ws.Range("H89").Value = 708.1818
$ws.Range("I89").Value = 656.8
$ws.Range("J89").Value = 751
$ws.Range("K89").Value = 5911.2
$ws.Range("L89").Value = 6759
$ws.Range("M89").Value = 16.80000000000018
$ws.Range("N89").Value = -18615

$ws.Range("H107").Value = 343.33334
$ws.Range("J107").Value = 376.1111
$ws.Range("L107").Value = 1128.3333
$ws.Range("N107").Value = -4968.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6290.9165
$ws.Range("I126").Value = 7413.6665
$ws.Range("K126").Value = 22240.9995
$ws.Range("M126").Value = -19770.9995

$ws.Range("H132").Value = 5370.3335
$ws.Range("I132").Value = 5692.9414
$ws.Range("K132").Value = 17078.8242
$ws.Range("M132").Value = -14548.8242

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3616.4783
$ws.Range("I7").Value = 3413.5
$ws.Range("J7").Value = 3932.2222
$ws.Range("K7").Value = 3413.5
$ws.Range("L7").Value = 3932.2222
$ws.Range("M7").Value = -3301.5
$ws.Range("N7").Value = -4156.2222

$ws.Range("H61").Value = 839.2
$ws.Range("I61").Value = 865.7778
$ws.Range("J61").Value = 600
$ws.Range("K61").Value = 865.7778
$ws.Range("L61").Value = 600
$ws.Range("M61").Value = -663.7778
$ws.Range("N61").Value = -1004

$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954

$ws.Range("H113").Value = 839.2
$ws.Range("I113").Value = 865.7778
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 865.7778
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = 1304.2222
$ws.Range("N113").Value = -4940

$ws.Range("H126").Value = 3616.4783
$ws.Range("I126").Value = 3413.5
$ws.Range("J126").Value = 3932.2222
$ws.Range("K126").Value = 10240.5
$ws.Range("L126").Value = 11796.6666
$ws.Range("M126").Value = -7770.5
$ws.Range("N126").Value = -16736.6666

$ws.Range("H132").Value = 2197
$ws.Range("I132").Value = 2197
$ws.Range("K132").Value = 6591
$ws.Range("M132").Value = -4061

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 298.5
$ws.Range("I107").Value = 298.5
$ws.Range("K107").Value = 895.5
$ws.Range("M107").Value = 1024.5

$ws.Range("H112").Value = 20387
$ws.Range("J112").Value = 20387
$ws.Range("L112").Value = 20387
$ws.Range("N112").Value = -23341

$ws.Range("H113").Value = 1158.6086
$ws.Range("I113").Value = 409.3889
$ws.Range("J113").Value = 3855.8
$ws.Range("K113").Value = 1228.1667
$ws.Range("L113").Value = 11567.4
$ws.Range("M113").Value = 941.8333
$ws.Range("N113").Value = -15907.4

$ws.Range("H136").Value = 4048.742
$ws.Range("I136").Value = 4432.909
$ws.Range("K136").Value = 13298.727
$ws.Range("M136").Value = -10748.727
